# "update scripts wuth new tpm"
#
# The NATMI ligand-receptor (Btla -> Tnfrsf14) table was rebuilt against a
# refreshed TPM expression matrix. The sending/ligand/receptor/target-cluster
# labels in columns A:D are unchanged, but:
#   - "MuSCs" no longer qualifies as a sending cluster under the new TPM
#     values, so the three rows where it was the sender (old rows 8-10) are
#     removed entirely, shrinking the table from 9 data rows to 6
#     (dimension A1:T10 -> A1:T7).
#   - Every remaining row's computed NATMI statistics (columns E:T) are
#     recalculated from the new TPM data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old "MuSCs" sending-cluster block (rows 8-10).
$ws.Rows("8:10").Delete()

# New per-row NATMI statistics (columns E through T) computed from the
# refreshed TPM values. Row/column labels (A:D) are untouched.
$cols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$newValues = @{
    2 = @(2, 0.6666666666666666, 0.3107469999999999, 0.9322409999999999,
          0.9278633407583023, 0.9278633407583025, 3, 1,
          3.558321333333333, 10.674964, 0.3039644761000113, 0.3039644761000113,
          1.105737679369333, 9.951639114323998, 0.2820374942660037, 0.2820374942660036)
    3 = @(2, 0.6666666666666666, 0.3107469999999999, 0.9322409999999999,
          0.9278633407583023, 0.9278633407583025, 3, 1,
          5.383140666666667, 16.149422, 0.4598470400038817, 0.4598470400038817,
          1.672794812744666, 15.055153314702, 0.4266752107758184, 0.4266752107758184)
    4 = @(2, 0.6666666666666666, 0.3107469999999999, 0.9322409999999999,
          0.9278633407583023, 0.9278633407583025, 3, 1,
          2.764910333333333, 8.294730999999999, 0.2361884838961071, 0.236188483896107,
          0.859187591352333, 7.732688322170998, 0.2191506357164804, 0.2191506357164804)
    5 = @(1, 0.3333333333333333, 0.024159, 0.072477,
          0.07213665924169768, 0.07213665924169768, 3, 1,
          3.558321333333333, 10.674964, 0.3039644761000113, 0.3039644761000113,
          0.085965485092, 0.773689365828, 0.02192698183400768, 0.02192698183400767)
    6 = @(1, 0.3333333333333333, 0.024159, 0.072477,
          0.07213665924169768, 0.07213665924169768, 3, 1,
          5.383140666666667, 16.149422, 0.4598470400038817, 0.4598470400038817,
          0.130051295366, 1.170461658294, 0.03317182922806334, 0.03317182922806333)
    7 = @(1, 0.3333333333333333, 0.024159, 0.072477,
          0.07213665924169768, 0.07213665924169768, 3, 1,
          2.764910333333333, 8.294730999999999, 0.2361884838961071, 0.236188483896107,
          0.066797468743, 0.6011772186869999, 0.01703784817962668, 0.01703784817962667)
}

foreach ($r in @(2,3,4,5,6,7)) {
    $vals = $newValues[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $vals[$i]
    }
}
